# Working-hours sheet: add a new work entry (row) before the Total row,
# shifting the Total row down and extending its SUM() ranges.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 5, pushing the old row 5 (Total) to row 6.
$ws.Rows("5:5").Insert() | Out-Null

# Seed the new row's formatting from the row above (an existing data row)
# so it picks up the same date / time / number styles.
$ws.Range("A2:F2").Copy() | Out-Null
$ws.Range("A5:F5").PasteSpecial(-4122) | Out-Null

# New work entry: 2023-12-07, 15:00 -> 21:00, rate 10/hr.
$ws.Range("A5").Value = 45267
$ws.Range("B5").Value = 0.625
$ws.Range("C5").Value = 0.875
$ws.Range("D5").Formula = "=(C5<B5)+C5-B5"
$ws.Range("E5").Value = 10
$ws.Range("F5").Formula = "=(D5*24)*E5"

# Re-apply the data-row formatting so the formula cells (which can pick up
# a borrowed number format from their precedents) end up styled the same
# as the rest of the data rows.
$ws.Range("A2:F2").Copy() | Out-Null
$ws.Range("A5:F5").PasteSpecial(-4122) | Out-Null

# The Total row, now on row 6, must sum across the new row too.
$ws.Range("D6").Formula = "=SUM(D2:D5)"
$ws.Range("F6").Formula = "=SUM(F2:F5)"

# Match the saved selection from the edit.
$ws.Range("F6").Select() | Out-Null
